$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Fill in the "15 min" entry that used to be a blank template row (15/16)
# ---------------------------------------------------------------------------
$sheet1.Range("A15").Value = 43901
$sheet1.Range("B15").Value = "15 min"
$sheet1.Range("C15").Value = "Explication"
$sheet1.Range("D15").Value = "ICT-431"
$sheet1.Range("E15").Value = "Problem avec git. Un merge as du etre fait car il y avais une différence entre git et github. "
$sheet1.Range("G15").Value = "Aide de M.Favre."
$sheet1.Range("K15").Value = 5
$sheet1.Rows.Item(15).RowHeight = 47.25

$sheet1.Range("A16").Value = 43901
$sheet1.Range("B16").Value = "15 min"

# ---------------------------------------------------------------------------
# 2. Update the current selection on "Journal de travail" (row 10 selected)
# ---------------------------------------------------------------------------
$sheet1.Range("A10:XFD10").Select()

# ---------------------------------------------------------------------------
# 3. Add the new "Journal de board" sheet right after "Journal de travail"
# ---------------------------------------------------------------------------
$sheet2 = $wb.Worksheets.Add($null, $sheet1)
$sheet2.Name = "Journal de board"

Write-Host "done"
